# Auto-generated edit script for Bmp7-Bmpr1a.xlsx (NATMI output)
# Commit: "Natmi following Dr Hou advice"
# Rewrites rows 2-6 (changed numbers/target clusters) and appends new rows 7-13
# for the "Neutro" sending cluster and rows where "sCs" is the sender across all
# target clusters (ECs, FAPs, M1, M2, Neutro, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20
# Row 2
$data[0,0] = "Neutro"
$data[0,1] = "Bmp7"
$data[0,2] = "Bmpr1a"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.1301303333333333
$data[0,7] = 0.390391
$data[0,8] = 0.4028730131292778
$data[0,9] = 0.5029891437391611
$data[0,10] = 2
$data[0,11] = 1
$data[0,12] = 4.344454
$data[0,13] = 8.688908
$data[0,14] = 0.07166328453363975
$data[0,15] = 0.05740743684517152
$data[0,16] = 0.5653452471713333
$data[0,17] = 3.392071483028
$data[0,18] = 0.02887120337080822
$data[0,19] = 0.0288753175030128
# Row 3
$data[1,0] = "Neutro"
$data[1,1] = "Bmp7"
$data[1,2] = "Bmpr1a"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.1301303333333333
$data[1,7] = 0.390391
$data[1,8] = 0.4028730131292778
$data[1,9] = 0.5029891437391611
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 24.18506433333333
$data[1,13] = 72.555193
$data[1,14] = 0.3989410744788757
$data[1,15] = 0.4793706711978917
$data[1,16] = 3.147210483384777
$data[1,17] = 28.324894350463
$data[1,18] = 0.1607225927363363
$data[1,19] = 0.2411182434394945
# Row 4
$data[2,0] = "Neutro"
$data[2,1] = "Bmp7"
$data[2,2] = "Bmpr1a"
$data[2,3] = "M1"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.1301303333333333
$data[2,7] = 0.390391
$data[2,8] = 0.4028730131292778
$data[2,9] = 0.5029891437391611
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.1434473333333333
$data[2,13] = 0.430342
$data[2,14] = 0.002366213812888463
$data[2,15] = 0.002843260762667162
$data[2,16] = 0.01866684930244444
$data[2,17] = 0.168001643722
$data[2,18] = 0.0009532836885064922
$data[2,19] = 0.00143012929644111
# Row 5
$data[3,0] = "Neutro"
$data[3,1] = "Bmp7"
$data[3,2] = "Bmpr1a"
$data[3,3] = "M2"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.1301303333333333
$data[3,7] = 0.390391
$data[3,8] = 0.4028730131292778
$data[3,9] = 0.5029891437391611
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 1.071697333333334
$data[3,13] = 3.215092
$data[3,14] = 0.01767802143436429
$data[3,15] = 0.0212420468649704
$data[3,16] = 0.1394603312191111
$data[3,17] = 1.255142980972
$data[3,18] = 0.007121997761426298
$data[3,19] = 0.01068451896387859
# Row 6
$data[4,0] = "Neutro"
$data[4,1] = "Bmp7"
$data[4,2] = "Bmpr1a"
$data[4,3] = "Neutro"
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.1301303333333333
$data[4,7] = 0.390391
$data[4,8] = 0.4028730131292778
$data[4,9] = 0.5029891437391611
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 4.708586333333334
$data[4,13] = 14.125759
$data[4,14] = 0.07766977441972553
$data[4,15] = 0.09332859982895587
$data[4,16] = 0.6127299090854444
$data[4,17] = 5.514569181769
$data[4,18] = 0.03129105604954613
$data[4,19] = 0.04694327251434133
# Row 7
$data[5,0] = "Neutro"
$data[5,1] = "Bmp7"
$data[5,2] = "Bmpr1a"
$data[5,3] = "sCs"
$data[5,4] = 1
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.1301303333333333
$data[5,7] = 0.390391
$data[5,8] = 0.4028730131292778
$data[5,9] = 0.5029891437391611
$data[5,10] = 2
$data[5,11] = 1
$data[5,12] = 26.1699
$data[5,13] = 52.3398
$data[5,14] = 0.4316816313205064
$data[5,15] = 0.3458079845003432
$data[5,16] = 3.4054978103
$data[5,17] = 20.4329868618
$data[5,18] = 0.1739128795226544
$data[5,19] = 0.1739376620219928
# Row 8
$data[6,0] = "sCs"
$data[6,1] = "Bmp7"
$data[6,2] = "Bmpr1a"
$data[6,3] = "ECs"
$data[6,4] = 2
$data[6,5] = 1
$data[6,6] = 0.1928755
$data[6,7] = 0.385751
$data[6,8] = 0.5971269868707222
$data[6,9] = 0.4970108562608389
$data[6,10] = 2
$data[6,11] = 1
$data[6,12] = 4.344454
$data[6,13] = 8.688908
$data[6,14] = 0.07166328453363975
$data[6,15] = 0.05740743684517152
$data[6,16] = 0.837938737477
$data[6,17] = 3.351754949908
$data[6,18] = 0.04279208116283153
$data[6,19] = 0.02853211934215873
# Row 9
$data[7,0] = "sCs"
$data[7,1] = "Bmp7"
$data[7,2] = "Bmpr1a"
$data[7,3] = "FAPs"
$data[7,4] = 2
$data[7,5] = 1
$data[7,6] = 0.1928755
$data[7,7] = 0.385751
$data[7,8] = 0.5971269868707222
$data[7,9] = 0.4970108562608389
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 24.18506433333333
$data[7,13] = 72.555193
$data[7,14] = 0.3989410744788757
$data[7,15] = 0.4793706711978917
$data[7,16] = 4.664706375823833
$data[7,17] = 27.988238254943
$data[7,18] = 0.2382184817425394
$data[7,19] = 0.2382524277583972
# Row 10
$data[8,0] = "sCs"
$data[8,1] = "Bmp7"
$data[8,2] = "Bmpr1a"
$data[8,3] = "M1"
$data[8,4] = 2
$data[8,5] = 1
$data[8,6] = 0.1928755
$data[8,7] = 0.385751
$data[8,8] = 0.5971269868707222
$data[8,9] = 0.4970108562608389
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.1434473333333333
$data[8,13] = 0.430342
$data[8,14] = 0.002366213812888463
$data[8,15] = 0.002843260762667162
$data[8,16] = 0.02766747614033334
$data[8,17] = 0.166004856842
$data[8,18] = 0.001412930124381971
$data[8,19] = 0.001413131466226052
# Row 11
$data[9,0] = "sCs"
$data[9,1] = "Bmp7"
$data[9,2] = "Bmpr1a"
$data[9,3] = "M2"
$data[9,4] = 2
$data[9,5] = 1
$data[9,6] = 0.1928755
$data[9,7] = 0.385751
$data[9,8] = 0.5971269868707222
$data[9,9] = 0.4970108562608389
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 1.071697333333334
$data[9,13] = 3.215092
$data[9,14] = 0.01767802143436429
$data[9,15] = 0.0212420468649704
$data[9,16] = 0.2067041590153334
$data[9,17] = 1.240224954092
$data[9,18] = 0.01055602367293799
$data[9,19] = 0.01055752790109181
# Row 12
$data[10,0] = "sCs"
$data[10,1] = "Bmp7"
$data[10,2] = "Bmpr1a"
$data[10,3] = "Neutro"
$data[10,4] = 2
$data[10,5] = 1
$data[10,6] = 0.1928755
$data[10,7] = 0.385751
$data[10,8] = 0.5971269868707222
$data[10,9] = 0.4970108562608389
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 4.708586333333334
$data[10,13] = 14.125759
$data[10,14] = 0.07766977441972553
$data[10,15] = 0.09332859982895587
$data[10,16] = 0.9081709433348334
$data[10,17] = 5.449025660009
$data[10,18] = 0.0463787183701794
$data[10,19] = 0.04638532731461454
# Row 13
$data[11,0] = "sCs"
$data[11,1] = "Bmp7"
$data[11,2] = "Bmpr1a"
$data[11,3] = "sCs"
$data[11,4] = 2
$data[11,5] = 1
$data[11,6] = 0.1928755
$data[11,7] = 0.385751
$data[11,8] = 0.5971269868707222
$data[11,9] = 0.4970108562608389
$data[11,10] = 2
$data[11,11] = 1
$data[11,12] = 26.1699
$data[11,13] = 52.3398
$data[11,14] = 0.4316816313205064
$data[11,15] = 0.3458079845003432
$data[11,16] = 5.047532547449999
$data[11,17] = 20.1901301898
$data[11,18] = 0.2577687517978519
$data[11,19] = 0.1718703224783505

# Write the full A2:T13 block in one shot
$ws.Range("A2:T13").Value = $data

